# Updated OPD billing test script
# Rewrites the "test" sheet's TestCaseName (A) / Plan (C) columns for the
# Reports section, inserts several new report test cases, appends three
# brand-new rows (TC022-TC024), and moves the " " marker from J14 to J17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# --- Column A (TestCaseName) updates for rows whose test script changed ---
$ws.Range("A7").Value  = "Reports\TC011UserCollectionReport.py"
$ws.Range("A8").Value  = "Reports\TC002TotalItemsBillReport.py"
$ws.Range("A9").Value  = "Reports\TC005IncomeSegregationReport.py"
$ws.Range("A10").Value = "Reports\TC013CancelBillReport.py"
$ws.Range("A11").Value = "Reports\TC014ReturnBillsReport.py"
$ws.Range("A12").Value = "Reports\TC008DiscountReport.py"
$ws.Range("A13").Value = "Reports\TC015EHSBillReport.py"
$ws.Range("A14").Value = "Laboratory\TC005GenerateLabReport.py"
$ws.Range("A15").Value = "Radiology\TC001GenerateUSGReport.py"
$ws.Range("A16").Value = "Dispensary\TC001CreateDispensarySale.py"
$ws.Range("A17").Value = "Pharmacy\Reports\TC002UserCollectionReport.py"
$ws.Range("A18").Value = "ADT\TC010AdmissionDischargeTransferToBePaid.py"
$ws.Range("A19").Value = "ADT\TC011AdmissionDischargeTransferNoDeposit.py"
$ws.Range("A20").Value = "Nursing\TC001WardBilling&AddVitals.py"
$ws.Range("A21").Value = "Inventory\TC001CreateGoodReceipt.py"
$ws.Range("A22").Value = "MedicalRecords\TC001createBirth&DeathCertificate.py"

# --- Column C (Plan) updates to match the re-ordered rows ---
$ws.Range("C7").Value  = "SmokeTest"
$ws.Range("C8").Value  = "SmokeTest"
$ws.Range("C9").Value  = "SmokeTest"
$ws.Range("C10").Value = "SmokeTest"
$ws.Range("C11").Value = "SmokeTest"
$ws.Range("C12").Value = "SanityTest"
$ws.Range("C13").Value = "SanityTest"
$ws.Range("C14").Value = "SmokeTest"
$ws.Range("C15").Value = "SmokeTest"
$ws.Range("C16").Value = "SmokeTest"
$ws.Range("C17").Value = "SmokeTest"

# --- Move the " " marker that lived in J14 down to J17 (row shifted) ---
$ws.Range("J14").ClearContents()
$ws.Range("J17").Value = " "

# --- Append three brand-new rows (TC022, TC023, TC024) using the same
#     formatting as the existing body rows ---
$ws.Range("A22:H22").Copy() | Out-Null
$ws.Range("A23:H25").PasteSpecial(-4122) | Out-Null

$ws.Range("A23").Value = "MedicalRecords\Reports\TC001HospitalServiceSummaryReport.py"
$ws.Range("B23").Value = "Norun"
$ws.Range("C23").Value = "SanityTest"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = "TC022"
$ws.Range("H23").Value = "LPH"

$ws.Range("A24").Value = "MedicalRecords\Reports\TC002InpatientMorbidityReport.py"
$ws.Range("B24").Value = "Norun"
$ws.Range("C24").Value = "SanityTest"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = "TC023"
$ws.Range("H24").Value = "LPH"

$ws.Range("A25").Value = "Vaccination\TC001RegisterVaccinationPatient.py"
$ws.Range("B25").Value = "Norun"
$ws.Range("C25").Value = "SanityTest"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "TC024"
$ws.Range("H25").Value = "LPH"

# --- Match the final selection recorded in the saved workbook ---
$ws.Range("B29").Select()
